$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two elevator trim-curve measurement rows in the table (originally
# rows 59 and 61) were reordered: the row that used to be first (A=1 ...)
# now comes after the one that used to be third (A=3 ...). Row 60, the
# measurement in between, is untouched. Rather than moving rows, swap the
# cell values between row 59 and row 61 so the sheet ends up with the new
# order while every other row stays exactly where it was.

# Row 59 now holds what used to be row 61's measurement.
$ws.Range("A59").Value = 3
$ws.Range("B59").Value = 0.024999999999999998
$ws.Range("C59").Value = 54
$ws.Range("D59").Value = 9660
$ws.Range("E59").Value = 150
$ws.Range("F59").Value = "6,2"
$ws.Range("G59").Value = "-0,2"
$ws.Range("H59").Value = "3,7"
$ws.Range("I59").Value = "-35,5"
$ws.Range("J59").Value = 446
$ws.Range("K59").Value = 488
$ws.Range("L59").Value = 622
$ws.Range("M59").Value = -7

# Row 61 now holds what used to be row 59's measurement.
$ws.Range("A61").Value = 1
$ws.Range("B61").Value = 0.022222222222222223
$ws.Range("C61").Value = 46
$ws.Range("D61").Value = 9210
$ws.Range("E61").Value = 170
$ws.Range("F61").Value = "4,6"
$ws.Range("G61").Value = "0,5"
$ws.Range("H61").Value = "3,7"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 454
$ws.Range("K61").Value = 495
$ws.Range("L61").Value = 591
$ws.Range("M61").Value = "-4,8"

# Row 58 was a blank spacer row whose cells still carried leftover
# formatting; clear it out so the row is truly empty.
$ws.Range("A58:M58").Clear()

# Leave the selection where the edit ended up, matching the saved view.
$null = $ws.Range("L64").Select()
